$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "AddStudy" label to "AddNewStudy"
$ws.Range("A1").Value = "AddNewStudy"

# Update the active selection to A2
$ws.Range("A2").Select()
